# Auto-generated edit script: updates cached market-derived profit
# figures (columns H-N) on each class sheet, per the scheduled-runner
# refresh described in the commit. Values are static numbers (no
# formulas live in these cells), so we just rewrite them directly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 324.66666
$ws.Range("I38").Value = 252.75
$ws.Range("K38").Value = 758.25
$ws.Range("M38").Value = -386.25
$ws.Range("H58").Value = 1156.4
$ws.Range("I58").Value = 1062.6666
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 3187.9998
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -3037.9998
$ws.Range("N58").Value = -6300
$ws.Range("H87").Value = 24128.322
$ws.Range("J87").Value = 24128.322
$ws.Range("L87").Value = 24128.322
$ws.Range("N87").Value = -26624.322
$ws.Range("H90").Value = 24128.322
$ws.Range("J90").Value = 24128.322
$ws.Range("L90").Value = 72384.966
$ws.Range("N90").Value = -84864.966
$ws.Range("H98").Value = 1683.3334
$ws.Range("I98").Value = 1525
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 1525
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -27
$ws.Range("N98").Value = -4996
$ws.Range("H122").Value = 1683.3334
$ws.Range("I122").Value = 1525
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4575
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2125
$ws.Range("N122").Value = -10900
$ws.Range("H125").Value = 2414.4
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2414.4
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 21729.6
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -26649.6
$ws.Range("H137").Value = 2824.7112
$ws.Range("I137").Value = 2621.7942
$ws.Range("J137").Value = 3451.9092
$ws.Range("K137").Value = 7865.382599999999
$ws.Range("L137").Value = 10355.7276
$ws.Range("M137").Value = -5315.382599999999
$ws.Range("N137").Value = -15455.7276
$ws.Range("H138").Value = 2429.9722
$ws.Range("J138").Value = 2156.0908
$ws.Range("L138").Value = 6468.2724
$ws.Range("N138").Value = -16748.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2850.6086
$ws.Range("I45").Value = 2080.7693
$ws.Range("J45").Value = 3851.4
$ws.Range("K45").Value = 2080.7693
$ws.Range("L45").Value = 3851.4
$ws.Range("M45").Value = -1703.7693
$ws.Range("N45").Value = -4605.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 652.1111
$ws.Range("I64").Value = 616.6667
$ws.Range("J64").Value = 723
$ws.Range("K64").Value = 616.6667
$ws.Range("L64").Value = 723
$ws.Range("M64").Value = -391.6667
$ws.Range("N64").Value = -1173
$ws.Range("H67").Value = 652.1111
$ws.Range("I67").Value = 616.6667
$ws.Range("J67").Value = 723
$ws.Range("K67").Value = 616.6667
$ws.Range("L67").Value = 723
$ws.Range("M67").Value = 163.3333
$ws.Range("N67").Value = -2283

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 5485.4136
$ws.Range("I31").Value = 1085.9584
$ws.Range("J31").Value = 26602.8
$ws.Range("K31").Value = 1085.9584
$ws.Range("L31").Value = 26602.8
$ws.Range("M31").Value = -790.9584
$ws.Range("N31").Value = -27192.8
$ws.Range("H34").Value = 5485.4136
$ws.Range("I34").Value = 1085.9584
$ws.Range("J34").Value = 26602.8
$ws.Range("K34").Value = 1085.9584
$ws.Range("L34").Value = 26602.8
$ws.Range("M34").Value = -883.9584
$ws.Range("N34").Value = -27006.8
$ws.Range("H114").Value = 69342
$ws.Range("J114").Value = 69342
$ws.Range("L114").Value = 69342
$ws.Range("N114").Value = -78020

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 180.2
$ws.Range("I8").Value = 180.2
$ws.Range("K8").Value = 540.5999999999999
$ws.Range("M8").Value = -401.5999999999999
$ws.Range("H39").Value = 1327.4667
$ws.Range("I39").Value = 407.4
$ws.Range("J39").Value = 1787.5
$ws.Range("K39").Value = 1222.2
$ws.Range("L39").Value = 5362.5
$ws.Range("M39").Value = -928.1999999999998
$ws.Range("N39").Value = -5950.5
$ws.Range("H110").Value = 12682.7
$ws.Range("J110").Value = 15971.429
$ws.Range("L110").Value = 47914.287
$ws.Range("N110").Value = -56094.287
$ws.Range("H113").Value = 968.5217
$ws.Range("I113").Value = 454.2857
$ws.Range("K113").Value = 1362.8571
$ws.Range("M113").Value = 807.1428999999998
$ws.Range("H122").Value = 7758.0713
$ws.Range("J122").Value = 50999
$ws.Range("L122").Value = 458991
$ws.Range("N122").Value = -463891
$ws.Range("H139").Value = 2538.6
$ws.Range("I139").Value = 1196.2727
$ws.Range("J139").Value = 3315.7368
$ws.Range("K139").Value = 3588.8181
$ws.Range("L139").Value = 9947.2104
$ws.Range("M139").Value = 1551.1819
$ws.Range("N139").Value = -20227.2104
$ws.Range("H141").Value = 6658.4165
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 6658.4165
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 19975.2495
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -30335.2495

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5177.706
$ws.Range("I70").Value = 5066.375
$ws.Range("J70").Value = 5336.75
$ws.Range("K70").Value = 5066.375
$ws.Range("L70").Value = 5336.75
$ws.Range("M70").Value = -4796.375
$ws.Range("N70").Value = -5876.75
$ws.Range("H73").Value = 5177.706
$ws.Range("I73").Value = 5066.375
$ws.Range("J73").Value = 5336.75
$ws.Range("K73").Value = 5066.375
$ws.Range("L73").Value = 5336.75
$ws.Range("M73").Value = -4130.375
$ws.Range("N73").Value = -7208.75
$ws.Range("H122").Value = 5322.759
$ws.Range("J122").Value = 5725.385
$ws.Range("L122").Value = 17176.155
$ws.Range("N122").Value = -22076.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 111114690
$ws.Range("I7").Value = 250001900
$ws.Range("K7").Value = 250001900
$ws.Range("M7").Value = -250001788
$ws.Range("H122").Value = 2734.5
$ws.Range("I122").Value = 2469
$ws.Range("K122").Value = 7407
$ws.Range("M122").Value = -4957
$ws.Range("H126").Value = 111114690
$ws.Range("I126").Value = 250001900
$ws.Range("K126").Value = 750005700
$ws.Range("M126").Value = -750003230
$ws.Range("H132").Value = 3619
$ws.Range("I132").Value = 2925.9473
$ws.Range("J132").Value = 5813.6665
$ws.Range("K132").Value = 8777.841899999999
$ws.Range("L132").Value = 17440.9995
$ws.Range("M132").Value = -6247.841899999999
$ws.Range("N132").Value = -22500.9995
$ws.Range("H136").Value = 2901.5
$ws.Range("I136").Value = 3318.077
$ws.Range("J136").Value = 1818.4
$ws.Range("K136").Value = 9954.231
$ws.Range("L136").Value = 5455.200000000001
$ws.Range("M136").Value = -7404.231
$ws.Range("N136").Value = -10555.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 684.3333
$ws.Range("J107").Value = 672.5
$ws.Range("L107").Value = 2017.5
$ws.Range("N107").Value = -5857.5
$ws.Range("H126").Value = 1303.1578
$ws.Range("I126").Value = 1181.9231
$ws.Range("K126").Value = 3545.7693
$ws.Range("M126").Value = -1075.7693
$ws.Range("H136").Value = 2207.3877
$ws.Range("I136").Value = 1868.2903
$ws.Range("K136").Value = 5604.8709

